# "Possibilité de sauvegarder la BDD sur un fichier Excel via le module
# 'openpyxl'." — a new order was recorded for the customer Fernandez/Gonzalo
# (seller "Flammarion" instead of "FNAC", with a real order_date timestamp
# instead of plain text) and a new seller "Coiffard" (Nantes, France) was
# appended to the Sellers sheet.

$wb = $excel.ActiveWorkbook

# --- Orders sheet -----------------------------------------------------
$wsOrders = $wb.Worksheets.Item("Orders")

# Row 2's order_date: stored as literal text before -> real Excel datetime now.
$wsOrders.Range("E2").Value = 44809.6475347222

# Row 3: seller changed from "FNAC" to "Flammarion", and its order_date also
# becomes a real datetime value instead of text.
$wsOrders.Range("C3").Value = "Flammarion"
$wsOrders.Range("E3").Value = 44809.6466666667

# --- Sellers sheet ------------------------------------------------------
$wsSellers = $wb.Worksheets.Item("Sellers")

# New seller row appended under the existing "Durance" entry.
$wsSellers.Range("A3").Value = "Coiffard"
$wsSellers.Range("B3").Value = "Nantes"
$wsSellers.Range("C3").Value = "France"

# --- Selection / active sheet bookkeeping -------------------------------
# Orders keeps its cursor but moves off C3 rather than E3, and loses focus
# to Sellers, which becomes the active tab with C3 selected there too.
$null = $wsOrders.Range("C3").Select()

$wsSellers.Activate()
$null = $wsSellers.Range("C3").Select()
